# "control de coma por punto decimal"
# Updates the fixed-point-iteration table (tabla_pf) with recomputed
# xn / fxn / E values (now rendered using a decimal point instead of a
# comma) and trims the table from 15 data rows (iterations 0-14) down to
# 12 data rows (iterations 0-11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2..13 (iterations 0..11). Column A (iteration index)
# is unchanged; columns B (xn), C (fxn) and D (E) get new figures.
$data = @(
  @("0",  "3.0",                "-65.0",                "1.000005"),
  @("1",  "20.6902860206768",   "-2862815909955.03",    "17.6902860206768"),
  @("2",  "34.2076502103495",   "-3.93602265055717e+20","13.5173641896727"),
  @("3",  "37.7271451659029",   "-5.17612423364357e+22","3.51949495555342"),
  @("4",  "38.4126590639731",   "-1.33882922202356e+23","0.685513898070205"),
  @("5",  "38.5387094787372",   "-1.59446538199507e+23","0.126050414764144"),
  @("6",  "38.5616422402785",   "-1.64597029342916e+23","0.0229327615412629"),
  @("7",  "38.56580640638",     "-1.65549956350725e+23","0.0041641661015461"),
  @("8",  "38.5665622763913",   "-1.65723520171509e+23","0.0007558700112753"),
  @("9",  "38.566699471454",    "-1.65755042579763e+23","0.000137195062706"),
  @("10", "38.5667243729146",   "-1.65760764667644e+23","2.49014605913089e-05"),
  @("11", "38.5667288926213",   "-1.65761803268812e+23","4.51970664983037e-06")
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $r = 2 + $i
  $rowVals = $data[$i]
  for ($c = 1; $c -le 4; $c++) {
    $cell = $ws.Cells.Item($r, $c)
    # Force text storage so values keep their literal "." decimal
    # representation instead of being reinterpreted as locale-formatted
    # numbers (the whole point of this commit).
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$c - 1]
  }
}

# The previous table had 15 data rows (iterations 0-14, sheet rows 2-16);
# the new table only has 12 (iterations 0-11, sheet rows 2-13). Remove the
# now-unused trailing rows so the used range / dimension shrinks to A1:D13.
$ws.Range("A14:D16").Delete()

Write-Host "tabla_pf updated"
